{"js": "// Update the answer table: replace the two-digit \u00f7 one-digit division\n// answers cell-by-cell (by position), matching the regenerated output.\n// The table has 20 rows x 5 columns; every 4th row (0, 4, 8, 12, 16)\n// holds the visible \"a\u00f7b=c, d\" answers, the rows in between are spacers.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, newValue] \u2014 0-based indices into the table grid.\nconst updates = [\n  [0, 0, \"89\u00f78=11, 1\"],\n  [0, 1, \"83\u00f76=13, 5\"],\n  [0, 2, \"47\u00f78=5, 7\"],\n  [0, 3, \"35\u00f78=4, 3\"],\n  [0, 4, \"15\u00f73=5, 0\"],\n\n  [4, 0, \"76\u00f77=10, 6\"],\n  [4, 1, \"54\u00f79=6, 0\"],\n  [4, 2, \"88\u00f73=29, 1\"],\n  [4, 3, \"31\u00f75=6, 1\"],\n  [4, 4, \"33\u00f78=4, 1\"],\n\n  // row 8, col 0 (\"78\u00f76=13, 0\") is unchanged in the target output.\n  [8, 1, \"59\u00f72=29, 1\"],\n  [8, 2, \"53\u00f79=5, 8\"],\n  [8, 3, \"72\u00f76=12, 0\"],\n  [8, 4, \"38\u00f74=9, 2\"],\n\n  [12, 0, \"18\u00f73=6, 0\"],\n  [12, 1, \"42\u00f78=5, 2\"],\n  [12, 2, \"32\u00f79=3, 5\"],\n  [12, 3, \"87\u00f75=17, 2\"],\n  [12, 4, \"81\u00f79=9, 0\"],\n\n  [16, 0, \"84\u00f77=12, 0\"],\n  [16, 1, \"20\u00f72=10, 0\"],\n  [16, 2, \"31\u00f76=5, 1\"],\n  [16, 3, \"10\u00f75=2, 0\"],\n  [16, 4, \"68\u00f76=11, 2\"],\n];\n\nfor (const [r, c, text] of updates) {\n  const cell = table.getCell(r, c);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the answer table: replace the two-digit \u00f7 one-digit division\n# answers cell-by-cell (by position), matching the regenerated output.\n# The table has 20 rows x 5 columns; every 4th row (1, 5, 9, 13, 17 in\n# 1-based COM indexing) holds the visible \"a\u00f7b=c, d\" answers, the rows\n# in between are spacers.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row, Column (1-based) -> new cell text.\n$updates = @(\n    @(1, 1, \"89\u00f78=11, 1\"),\n    @(1, 2, \"83\u00f76=13, 5\"),\n    @(1, 3, \"47\u00f78=5, 7\"),\n    @(1, 4, \"35\u00f78=4, 3\"),\n    @(1, 5, \"15\u00f73=5, 0\"),\n\n    @(5, 1, \"76\u00f77=10, 6\"),\n    @(5, 2, \"54\u00f79=6, 0\"),\n    @(5, 3, \"88\u00f73=29, 1\"),\n    @(5, 4, \"31\u00f75=6, 1\"),\n    @(5, 5, \"33\u00f78=4, 1\"),\n\n    # Row 9, column 1 (\"78\u00f76=13, 0\") is unchanged in the target output.\n    @(9, 2, \"59\u00f72=29, 1\"),\n    @(9, 3, \"53\u00f79=5, 8\"),\n    @(9, 4, \"72\u00f76=12, 0\"),\n    @(9, 5, \"38\u00f74=9, 2\"),\n\n    @(13, 1, \"18\u00f73=6, 0\"),\n    @(13, 2, \"42\u00f78=5, 2\"),\n    @(13, 3, \"32\u00f79=3, 5\"),\n    @(13, 4, \"87\u00f75=17, 2\"),\n    @(13, 5, \"81\u00f79=9, 0\"),\n\n    @(17, 1, \"84\u00f77=12, 0\"),\n    @(17, 2, \"20\u00f72=10, 0\"),\n    @(17, 3, \"31\u00f76=5, 1\"),\n    @(17, 4, \"10\u00f75=2, 0\"),\n    @(17, 5, \"68\u00f76=11, 2\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $text\n}\n"}
